# Recursive_COT_Sandbox/Results.xlsx -- "Document Updates & New Graphing Updates"
#
# The five small "scorecard" blocks (title + header + two AVERAGE formulas,
# living in columns J:K) that used to sit right under each model's own
# 15-row data table are being consolidated into one contiguous stack of
# blocks directly under the very first table (rows 5-7), i.e.:
#
#   old J20:K22  (GPT-4o CoT Enhanced)      -> new J8:K10   (brand new copy)
#   old J35:K37  (GPT-o1-preview)           -> new J11:K13
#   old J50:K52  (Claude-3.5-Sonnet)        -> new J14:K16
#   old J65:K67  (Claude-3-Opus)            -> new J17:K19
#   old J80:K82  (Claude-3-Opus CoT)        -> new J20:K22  (overwrites the
#                                               block that used to live here)
#
# The four vacated source blocks (35:37, 50:52, 65:67, 80:82) are removed
# entirely (content + formatting + merge), shrinking those rows back down
# to the A:H data only.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Capture every old value/formula we still need, before any of the
#    destination cells (which overlap some of the source cells, e.g.
#    J20:K22 is both a source's neighbour and the final destination of
#    the 80:82 block) get overwritten.
# ---------------------------------------------------------------------
$blocks = @{}
foreach ($old in 20, 35, 50, 65, 80) {
    $blocks[$old] = @{
        title = $ws.Range("J$old").Formula
        hdrJ  = $ws.Range("J$($old+1)").Formula
        hdrK  = $ws.Range("K$($old+1)").Formula
        fJ    = $ws.Range("J$($old+2)").Formula
        fK    = $ws.Range("K$($old+2)").Formula
    }
}

# ---------------------------------------------------------------------
# 2. Helper: stamp a 3-row block (merged title / header / formula row)
#    at a destination starting row, formatting it like the untouched
#    template block at J5:K7, then filling in the supplied content.
# ---------------------------------------------------------------------
function Write-Block($destStart, $content, $formatDest) {
    $tRow = $destStart
    $hRow = $destStart + 1
    $fRow = $destStart + 2

    if ($formatDest) {
        $ws.Range("J$tRow`:K$tRow").Merge()

        $ws.Range("J5").Copy()
        $ws.Range("J$tRow").PasteSpecial(-4122)
        $ws.Range("K5").Copy()
        $ws.Range("K$tRow").PasteSpecial(-4122)

        $ws.Range("J6").Copy()
        $ws.Range("J$hRow").PasteSpecial(-4122)
        $ws.Range("K6").Copy()
        $ws.Range("K$hRow").PasteSpecial(-4122)

        $ws.Range("J7").Copy()
        $ws.Range("J$fRow").PasteSpecial(-4122)
        $ws.Range("K7").Copy()
        $ws.Range("K$fRow").PasteSpecial(-4122)
    }

    $ws.Range("J$tRow").Formula = $content.title
    $ws.Range("J$hRow").Formula = $content.hdrJ
    $ws.Range("K$hRow").Formula = $content.hdrK
    $ws.Range("J$fRow").Formula = $content.fJ
    $ws.Range("K$fRow").Formula = $content.fK
}

# ---------------------------------------------------------------------
# 3. Create the four brand-new blocks (8:10, 11:13, 14:16, 17:19).
# ---------------------------------------------------------------------
Write-Block 8  $blocks[20] $true
Write-Block 11 $blocks[35] $true
Write-Block 14 $blocks[50] $true
Write-Block 17 $blocks[65] $true

# ---------------------------------------------------------------------
# 4. Overwrite the existing J20:K22 block in place with the old 80:82
#    content (same formatting already in place, merge already exists).
# ---------------------------------------------------------------------
Write-Block 20 $blocks[80] $false

# ---------------------------------------------------------------------
# 5. Remove the four vacated source blocks entirely.
# ---------------------------------------------------------------------
foreach ($old in 35, 50, 65, 80) {
    $ws.Range("J$old`:K$old").UnMerge()
    $ws.Range("J$old`:K$($old+2)").Clear()
}

# ---------------------------------------------------------------------
# 6. Cosmetic: restore the active-cell selection recorded in the file.
# ---------------------------------------------------------------------
$ws.Range("K35").Select()
